# "Generate Report for Archive"
#
# The handoff/localization status moved on from "Ready for handoff" to
# "In Translation" for both in-flight files, on every sheet that surfaces
# the Status column (the Overview rollup as well as the per-locale
# zh-cn / de-de detail sheets). Excel re-flows the Status column a bit
# narrower afterwards to match the shorter text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn (col E) and de-de (col F) status columns ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F3").Value = "In Translation"
$wsOverview.Range("E:F").ColumnWidth = 12.5

# --- zh-cn sheet: Status column (col C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C3").Value = "In Translation"
$wsZhCn.Range("C:C").ColumnWidth = 12.5

# --- de-de sheet: Status column (col C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C3").Value = "In Translation"
$wsDeDe.Range("C:C").ColumnWidth = 12.5
